# Update(io): get delete update
# Append a new data row (row 4) to Sheet1, mirroring the existing rows'
# layout: most columns are blank for this record, with a handful of
# numeric counts filled in (columns A, F, O, T and AD).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 33   # column AG
$newRow  = 4

# Clear/blank out the whole new row first so every column in A4:AG4 gets an
# explicit (empty) entry, matching the shape of rows 2-3 above it.
for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item($newRow, $col).Value = ""
}

# Now fill in the handful of non-blank numeric values for this row.
$ws.Cells.Item($newRow, 1).Value  = 3    # A4
$ws.Cells.Item($newRow, 6).Value  = 5    # F4
$ws.Cells.Item($newRow, 15).Value = 3    # O4
$ws.Cells.Item($newRow, 20).Value = 5    # T4
$ws.Cells.Item($newRow, 30).Value = 13   # AD4
